# feat: add 2022-Q1 data
#
# The workbook has per-quarter fund-holding sheets (2020-Q4 .. 2021-Q3) plus
# a trailing "总计" (totals) roll-up sheet. This change inserts a new
# "2022-Q1" fund-holding sheet right before "总计", and updates "总计" with
# a new first data row summarizing 2022-Q1 (2 funds, 0 亿元 held), shifting
# the previously-existing rows down by one.
#
# Strategy:
#  - Rename the existing "总计" sheet (sheetId=5) to "2022-Q1" in place, so
#    it keeps all of its sheetPr/pageMargins/sheetView/style plumbing, then
#    fill it with the new fund-holding data (columns A-H).
#  - Duplicate that sheet (.Copy) right after itself to get a new sheet
#    (sheetId=6) with identical structural formatting, rename the copy back
#    to "总计", and refill it with the shifted totals table (columns A-D).
#  - Restore the originally-active sheet/tab at the end.

$wb = $excel.ActiveWorkbook

# A pristine, never-written-to cell carrying the bold/centered/thin-border
# "s=2" header/index style, used below as a copy/paste-format source. It
# MUST be a cell this script never overwrites: PasteSpecial(Formats) reads
# the source at paste time (not at Copy() time), so using a cell as its own
# upcoming style source — after stamping a transient NumberFormat="@" text
# format onto it — would re-capture that transient state instead of the
# original look.
$styleSrc = $wb.Worksheets.Item("2021-Q3").Range("C1")
$idxStyleSrc = $wb.Worksheets.Item("2021-Q3").Range("A2")

# Helper: write a value into a cell as literal TEXT (shared-string), even if
# it looks like a number (e.g. "0.34", "009327") — matches the source data's
# t="inlineStr" cells, which must not be coerced into floats / lose leading
# zeros. Setting NumberFormat="@" before the assignment is what keeps Excel
# from re-interpreting the string as a number.
#
# That NumberFormat="@" stamp leaves a stray text-format style on the cell,
# so afterwards we reapply the "real" desired look:
#   - $styleSource omitted -> reset to the sheet default ("Normal"), for
#     plain unstyled data cells.
#   - $styleSource given   -> copy/paste just the *formats* from that
#     (already correctly-styled) cell on top, for header/index cells that
#     must keep a non-default style (e.g. the bold/centered/bordered "s=2"
#     look). PasteSpecial(Formats) only touches style, not value/type, so
#     the text we just assigned survives untouched.
function Set-TextValue($range, [string]$value, $styleSource = $null) {
    if ($styleSource -ne $null) {
        $styleSource.Copy()
    }
    $range.NumberFormat = "@"
    $range.Value = $value
    if ($styleSource -ne $null) {
        $range.PasteSpecial(-4122)   # xlPasteFormats
    } else {
        $range.Style = "Normal"
    }
}

# Helper: write the little "s=2"-styled row-index number that column A
# carries on every data sheet (0, 1, 2, ...) — plain numbers never get
# mis-coerced, so this only has to (re)stamp the style.
function Set-IndexValue($range, $value) {
    $idxStyleSrc.Copy()
    $range.PasteSpecial(-4122)   # xlPasteFormats
    $range.Value = $value
}

# ---------------------------------------------------------------------
# 1) "总计" -> "2022-Q1": new per-fund holdings sheet
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Extend the bold/centered/bordered header style (currently only on B1:D1)
# across the new E1:H1 header cells before we touch any values.
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)   # xlPasteFormats

# Wipe all prior "总计" values (old 2021-Q3/Q2/Q1/2020-Q4 rows lived in
# A2:D5) and then drop the now-unused rows 4:5 — the new sheet only needs
# rows 1-3. A2/A3 already carry "s=2" and only ever get plain-number
# values, so they need no special handling — ClearContents() leaves their
# style alone.
$q1.Range("A1:H10").ClearContents()
$q1.Rows("4:5").Delete()

# Header row
Set-TextValue $q1.Range("B1") "基金代码" $styleSrc
Set-TextValue $q1.Range("C1") "基金名称" $styleSrc
Set-TextValue $q1.Range("D1") "基金规模" $styleSrc
Set-TextValue $q1.Range("E1") "股票总仓位" $styleSrc
Set-TextValue $q1.Range("F1") "仓位占比" $styleSrc
Set-TextValue $q1.Range("G1") "持有市值(亿元)" $styleSrc
Set-TextValue $q1.Range("H1") "仓位排名" $styleSrc

# Row 2: 009327 东兴兴晟混合A
Set-IndexValue $q1.Range("A2") 0
Set-TextValue $q1.Range("B2") "009327"
Set-TextValue $q1.Range("C2") "东兴兴晟混合A"
Set-TextValue $q1.Range("D2") "0.34"
Set-TextValue $q1.Range("E2") "79.83"
Set-TextValue $q1.Range("F2") "0.97"
Set-TextValue $q1.Range("G2") "0.0033"
$q1.Range("H2").Value = 6

# Row 3: 009328 东兴兴晟混合C
Set-IndexValue $q1.Range("A3") 1
Set-TextValue $q1.Range("B3") "009328"
Set-TextValue $q1.Range("C3") "东兴兴晟混合C"
Set-TextValue $q1.Range("D3") "0.08"
Set-TextValue $q1.Range("E3") "79.83"
Set-TextValue $q1.Range("F3") "0.97"
Set-TextValue $q1.Range("G3") "0.0008"
$q1.Range("H3").Value = 6

# ---------------------------------------------------------------------
# 2) New "总计" sheet (duplicate of "2022-Q1" for identical formatting)
# ---------------------------------------------------------------------
$q1.Copy($null, $q1)
$total = $wb.Worksheets.Item("2022-Q1 (2)")
$total.Name = "总计"

$total.Range("A1:H10").ClearContents()
$total.Rows("4:8").Delete()
$total.Columns("E:H").Delete()

# Header row
Set-TextValue $total.Range("B1") "日期" $styleSrc
Set-TextValue $total.Range("C1") "持有数量(只)" $styleSrc
Set-TextValue $total.Range("D1") "持有市值(亿元)" $styleSrc

# Row 2: new 2022-Q1 summary
Set-IndexValue $total.Range("A2") 0
Set-TextValue $total.Range("B2") "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0

# Row 3: was row 2 (2021-Q3)
Set-IndexValue $total.Range("A3") 1
Set-TextValue $total.Range("B3") "2021-Q3"
$total.Range("C3").Value = 5
$total.Range("D3").Value = 0.09

# Row 4: was row 3 (2021-Q2)
Set-IndexValue $total.Range("A4") 2
Set-TextValue $total.Range("B4") "2021-Q2"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.03

# Row 5: was row 4 (2021-Q1)
Set-IndexValue $total.Range("A5") 3
Set-TextValue $total.Range("B5") "2021-Q1"
$total.Range("C5").Value = 7
$total.Range("D5").Value = 3.72

# Row 6: was row 5 (2020-Q4)
Set-IndexValue $total.Range("A6") 4
Set-TextValue $total.Range("B6") "2020-Q4"
$total.Range("C6").Value = 1
$total.Range("D6").Value = 0.04

# ---------------------------------------------------------------------
# 3) Restore original tab selection (unchanged by this edit)
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
